$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 54783.023
$ws.Range("J17").Value = 54783.023
$ws.Range("L17").Value = 164349.069
$ws.Range("N17").Value = -164685.069

# Row 64
$ws.Range("H64").Value = 3450
$ws.Range("I64").Value = 3442.8572
$ws.Range("J64").Value = 3500
$ws.Range("K64").Value = 3442.8572
$ws.Range("L64").Value = 3500
$ws.Range("M64").Value = -3194.8572
$ws.Range("N64").Value = -3996

# Row 67
$ws.Range("H67").Value = 3450
$ws.Range("I67").Value = 3442.8572
$ws.Range("J67").Value = 3500
$ws.Range("K67").Value = 3442.8572
$ws.Range("L67").Value = 3500
$ws.Range("M67").Value = -2584.8572
$ws.Range("N67").Value = -5216

# Row 74
$ws.Range("H74").Value = 3836.1333
$ws.Range("I74").Value = 3420.2222
$ws.Range("K74").Value = 3420.2222
$ws.Range("M74").Value = -2484.2222

# Row 76
$ws.Range("H76").Value = 2712.7646
$ws.Range("I76").Value = 2657
$ws.Range("K76").Value = 2657
$ws.Range("M76").Value = -2342

# Row 77
$ws.Range("H77").Value = 3836.1333
$ws.Range("I77").Value = 3420.2222
$ws.Range("K77").Value = 17101.111
$ws.Range("M77").Value = -12421.111

# Row 79
$ws.Range("H79").Value = 2712.7646
$ws.Range("I79").Value = 2657
$ws.Range("K79").Value = 2657
$ws.Range("M79").Value = -1565

# Row 87
$ws.Range("H87").Value = 26655.625
$ws.Range("J87").Value = 26655.625
$ws.Range("L87").Value = 26655.625
$ws.Range("N87").Value = -29151.625

# Row 90
$ws.Range("H90").Value = 26655.625
$ws.Range("J90").Value = 26655.625
$ws.Range("L90").Value = 79966.875
$ws.Range("N90").Value = -92446.875

# Row 106
$ws.Range("H106").Value = 2786.2856
$ws.Range("I106").Value = 1901.6666
$ws.Range("K106").Value = 1901.6666
$ws.Range("M106").Value = -1270.6666

# Row 125
$ws.Range("H125").Value = 1686.0526
$ws.Range("I125").Value = 1379.9
$ws.Range("J125").Value = 2026.2222
$ws.Range("K125").Value = 12419.1
$ws.Range("L125").Value = 18235.9998
$ws.Range("M125").Value = -9959.1
$ws.Range("N125").Value = -23155.9998


$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1467.127
$ws.Range("I45").Value = 1113.2084
$ws.Range("J45").Value = 2599.6667
$ws.Range("K45").Value = 1113.2084
$ws.Range("L45").Value = 2599.6667
$ws.Range("M45").Value = -736.2084
$ws.Range("N45").Value = -3353.6667

# Row 61
$ws.Range("H61").Value = 3233.3
$ws.Range("I61").Value = 1047.5714
$ws.Range("J61").Value = 8333.333000000001
$ws.Range("K61").Value = 1047.5714
$ws.Range("L61").Value = 8333.333000000001
$ws.Range("M61").Value = -835.5714
$ws.Range("N61").Value = -8757.333000000001

# Row 74
$ws.Range("H74").Value = 1431.3334
$ws.Range("I74").Value = 1329.6
$ws.Range("K74").Value = 1329.6
$ws.Range("M74").Value = -455.5999999999999

# Row 77
$ws.Range("H77").Value = 1431.3334
$ws.Range("I77").Value = 1329.6
$ws.Range("K77").Value = 6648
$ws.Range("M77").Value = -2280

# Row 102
$ws.Range("H102").Value = 4750
$ws.Range("I102").Value = 3000
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -1378
$ws.Range("N102").Value = -13244

# Row 136
$ws.Range("H136").Value = 3233.3
$ws.Range("I136").Value = 1047.5714
$ws.Range("J136").Value = 8333.333000000001
$ws.Range("K136").Value = 3142.7142
$ws.Range("L136").Value = 24999.999
$ws.Range("M136").Value = -592.7142000000003
$ws.Range("N136").Value = -30099.999


$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1654.6562
$ws.Range("I105").Value = 1536.5
$ws.Range("K105").Value = 1536.5
$ws.Range("M105").Value = 210.5


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3551.8462
$ws.Range("I31").Value = 2584.1304
$ws.Range("J31").Value = 4942.9375
$ws.Range("K31").Value = 2584.1304
$ws.Range("L31").Value = 4942.9375
$ws.Range("M31").Value = -2289.1304
$ws.Range("N31").Value = -5532.9375

# Row 34
$ws.Range("H34").Value = 3551.8462
$ws.Range("I34").Value = 2584.1304
$ws.Range("J34").Value = 4942.9375
$ws.Range("K34").Value = 2584.1304
$ws.Range("L34").Value = 4942.9375
$ws.Range("M34").Value = -2382.1304
$ws.Range("N34").Value = -5346.9375

# Row 97
$ws.Range("H97").Value = 49380.953
$ws.Range("J97").Value = 49380.953
$ws.Range("L97").Value = 49380.953
$ws.Range("N97").Value = -51362.953

# Row 99
$ws.Range("H99").Value = 3000
$ws.Range("J99").Value = 3500
$ws.Range("L99").Value = 3500
$ws.Range("N99").Value = -6496

# Row 122
$ws.Range("H122").Value = 3570.923
$ws.Range("I122").Value = 2765.6365
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 8296.9095
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -5846.9095
$ws.Range("N122").Value = -28900

# Row 126
$ws.Range("H126").Value = 3000
$ws.Range("J126").Value = 3500
$ws.Range("L126").Value = 10500
$ws.Range("N126").Value = -15440


$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 2806.25
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 2806.25
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 8418.75
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -10414.75

# Row 78
$ws.Range("H78").Value = 2806.25
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 2806.25
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 25256.25
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -35240.25

# Row 131
$ws.Range("H131").Value = 1347.127
$ws.Range("I131").Value = 2476
$ws.Range("J131").Value = 1134.1321
$ws.Range("K131").Value = 7428
$ws.Range("L131").Value = 3402.3963
$ws.Range("M131").Value = -2388
$ws.Range("N131").Value = -13482.3963


$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4448.8335
$ws.Range("I80").Value = 3750
$ws.Range("J80").Value = 4798.25
$ws.Range("K80").Value = 3750
$ws.Range("L80").Value = 4798.25
$ws.Range("M80").Value = -2752
$ws.Range("N80").Value = -6794.25

# Row 83
$ws.Range("H83").Value = 4448.8335
$ws.Range("I83").Value = 3750
$ws.Range("J83").Value = 4798.25
$ws.Range("K83").Value = 18750
$ws.Range("L83").Value = 23991.25
$ws.Range("M83").Value = -13758
$ws.Range("N83").Value = -33975.25


$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3847845.2
$ws.Range("I7").Value = 6667897
$ws.Range("J7").Value = 2320.3635
$ws.Range("K7").Value = 6667897
$ws.Range("L7").Value = 2320.3635
$ws.Range("M7").Value = -6667785
$ws.Range("N7").Value = -2544.3635

# Row 22
$ws.Range("H22").Value = 90910360
$ws.Range("I22").Value = 125000410
$ws.Range("J22").Value = 3566.6667
$ws.Range("K22").Value = 125000410
$ws.Range("L22").Value = 3566.6667
$ws.Range("M22").Value = -125000115
$ws.Range("N22").Value = -4156.6667

# Row 27
$ws.Range("H27").Value = 90910360
$ws.Range("I27").Value = 125000410
$ws.Range("J27").Value = 3566.6667
$ws.Range("K27").Value = 125000410
$ws.Range("L27").Value = 3566.6667
$ws.Range("M27").Value = -125000303
$ws.Range("N27").Value = -3780.6667

# Row 40
$ws.Range("H40").Value = 2280.5
$ws.Range("I40").Value = 1540.375
$ws.Range("J40").Value = 3267.3333
$ws.Range("K40").Value = 1540.375
$ws.Range("L40").Value = 3267.3333
$ws.Range("M40").Value = -1404.375
$ws.Range("N40").Value = -3539.3333

# Row 68
$ws.Range("H68").Value = 1863.4615
$ws.Range("I68").Value = 1030.9524
$ws.Range("J68").Value = 5360
$ws.Range("K68").Value = 1030.9524
$ws.Range("L68").Value = 5360
$ws.Range("M68").Value = -281.9523999999999
$ws.Range("N68").Value = -6858

# Row 71
$ws.Range("H71").Value = 1863.4615
$ws.Range("I71").Value = 1030.9524
$ws.Range("J71").Value = 5360
$ws.Range("K71").Value = 5154.762
$ws.Range("L71").Value = 26800
$ws.Range("M71").Value = -1410.762
$ws.Range("N71").Value = -34288

# Row 100
$ws.Range("H100").Value = 2481.3
$ws.Range("I100").Value = 1413.25
$ws.Range("J100").Value = 3193.3333
$ws.Range("K100").Value = 1413.25
$ws.Range("L100").Value = 3193.3333
$ws.Range("M100").Value = -872.25
$ws.Range("N100").Value = -4275.3333

# Row 122
$ws.Range("H122").Value = 3121.1667
$ws.Range("I122").Value = 2572.6667
$ws.Range("J122").Value = 4766.6665
$ws.Range("K122").Value = 7718.000100000001
$ws.Range("L122").Value = 14299.9995
$ws.Range("M122").Value = -5268.000100000001
$ws.Range("N122").Value = -19199.9995

# Row 126
$ws.Range("H126").Value = 3847845.2
$ws.Range("I126").Value = 6667897
$ws.Range("J126").Value = 2320.3635
$ws.Range("K126").Value = 20003691
$ws.Range("L126").Value = 6961.0905
$ws.Range("M126").Value = -20001221
$ws.Range("N126").Value = -11901.0905

# Row 132
$ws.Range("H132").Value = 2512.8286
$ws.Range("I132").Value = 1613.4231
$ws.Range("K132").Value = 4840.2693
$ws.Range("M132").Value = -2310.2693


$ws = $wb.Worksheets.Item("WVR")
# Row 93
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -39992

# Row 96
$ws.Range("H96").Value = 1000
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 373
$ws.Range("N96").ClearContents()

# Row 122
$ws.Range("H122").Value = 347001.88
$ws.Range("I122").Value = 456099.78
$ws.Range("J122").Value = 4122.7144
$ws.Range("K122").Value = 1368299.34
$ws.Range("L122").Value = 12368.1432
$ws.Range("M122").Value = -1365849.34
$ws.Range("N122").Value = -17268.1432

# Row 126
$ws.Range("H126").Value = 3573727.2
$ws.Range("J126").Value = 12503811
$ws.Range("L126").Value = 37511433
$ws.Range("N126").Value = -37516373

